# Update localization status report for handoff generation
$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Status text changes: "In Translation" -> "Ready for handoff"
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# Latest HO Xliff Generate Date / Latest Handoff Datetime updates
$wsOverview.Range("G2").Value = "2016-09-03 17:06:45"
$wsDeDe.Range("H2").Value = "2016-09-03 17:06:45"
$wsZhCn.Range("H2").Value = "2016-09-03 17:06:41"

# Column width adjustments to fit the new, longer status text
$wsOverview.Range("E:E").ColumnWidth = 17.2159881591797
$wsOverview.Range("F:F").ColumnWidth = 17.2159881591797
$wsZhCn.Range("C:C").ColumnWidth = 17.2159881591797
$wsDeDe.Range("C:C").ColumnWidth = 17.2159881591797
